# B6-PowerPoint.pptx edit
#
# 1) Re-style the three data tables (slides 14, 15, 16) from the bare
#    "Table_0" style ({EA9421E6-EF6A-463F-95D5-A2284107B07A}) to the
#    standard gallery table style {134FD33B-7793-4378-A453-7D61CA0ADBBB}.
# 2) Swap the deck's design from the "Integral" theme over to the
#    default "Office Theme" (the theme that this file already carries
#    for its Notes Master), via the Design/SlideMaster COM surface.

$p = $ppt.ActivePresentation

$targetStyleId = "{134FD33B-7793-4378-A453-7D61CA0ADBBB}"

for ($i = 1; $i -le $p.Slides.Count; $i++) {
    $slide = $p.Slides.Item($i)
    for ($j = 1; $j -le $slide.Shapes.Count; $j++) {
        $shape = $slide.Shapes.Item($j)
        if ($shape.HasTable) {
            $shape.Table.ApplyStyle($targetStyleId)
        }
    }
}

# Apply the "Office Theme" design (the theme already embedded in this
# deck as the Notes Master's theme) to the slide master/design so the
# whole deck switches away from the Google-Slides-exported "Integral"
# theme.
$design = $p.Designs.Item(1)
$design.SlideMaster.ApplyTheme("Office Theme")
